# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
# Both sheets contain duplicated data, and both need the identical updates:
#   F17: 11   -> 12
#   F24: 3338 -> 3339
#   F25: 402  -> 403
#   F31: 1047 -> 1048

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F17").Value = 12
    $ws.Range("F24").Value = 3339
    $ws.Range("F25").Value = 403
    $ws.Range("F31").Value = 1048
}
